$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.054625478817785
$ws.Range("D2").Value = 1.06131217168353
$ws.Range("E2").Value = 1.060971200955186
$ws.Range("F2").Value = 1.071224200933356
$ws.Range("I2").Value = 1.041037728062605
$ws.Range("J2").Value = 1.059636593932684
$ws.Range("K2").Value = 1.064036574498938
$ws.Range("L2").Value = 1.063696530457602
$ws.Range("M2").Value = 1.073921943765801
$ws.Range("N2").Value = 1.023606697839558

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.055744417997425
$ws.Range("D3").Value = 1.062310261770003
$ws.Range("E3").Value = 1.061958125602051
$ws.Range("F3").Value = 1.072274699234095
$ws.Range("I3").Value = 1.041221652867747
$ws.Range("J3").Value = 1.060406068437119
$ws.Range("K3").Value = 1.064848798435474
$ws.Range("L3").Value = 1.0644975499802
$ws.Range("M3").Value = 1.07478837634353
$ws.Range("N3").Value = 1.023867585178404

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056468637870808
$ws.Range("D4").Value = 1.0629565452467
$ws.Range("E4").Value = 1.062597228101849
$ws.Range("F4").Value = 1.072954982001178
$ws.Range("I4").Value = 1.041339243277458
$ws.Range("J4").Value = 1.060903583739199
$ws.Range("K4").Value = 1.065374198063052
$ws.Range("L4").Value = 1.065015741048275
$ws.Range("M4").Value = 1.075348948594408
$ws.Range("N4").Value = 1.024036143672327

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056773146039486
$ws.Range("D5").Value = 1.063228350833391
$ws.Range("E5").Value = 1.062866025262084
$ws.Range("F5").Value = 1.073241101871123
$ws.Range("I5").Value = 1.041388337825004
$ws.Range("J5").Value = 1.061112646561037
$ws.Range("K5").Value = 1.065595036692842
$ws.Range("L5").Value = 1.06523355910598
$ws.Range("M5").Value = 1.075584596735301
$ws.Range("N5").Value = 1.024106944935762

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056824277015386
$ws.Range("D6").Value = 1.063273994512523
$ws.Range("E6").Value = 1.062911164425409
$ws.Range("F6").Value = 1.073289150211569
$ws.Range("I6").Value = 1.041396561044595
$ws.Range("J6").Value = 1.061147743693228
$ws.Range("K6").Value = 1.065632114155433
$ws.Range("L6").Value = 1.065270129981676
$ws.Range("M6").Value = 1.075624162123917
$ws.Range("N6").Value = 1.024118829214149

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056472706542188
$ws.Range("D7").Value = 1.062960176702215
$ws.Range("E7").Value = 1.062600819317578
$ws.Range("F7").Value = 1.072958804642412
$ws.Range("I7").Value = 1.041339900619248
$ws.Range("J7").Value = 1.060906377610962
$ws.Range("K7").Value = 1.065377149074612
$ws.Range("L7").Value = 1.065018651659899
$ws.Range("M7").Value = 1.075352097401297
$ws.Range("N7").Value = 1.024037089961072

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.055003590107941
$ws.Range("D8").Value = 1.061649387035491
$ws.Range("E8").Value = 1.061304633804516
$ws.Range("F8").Value = 1.071579109729429
$ws.Range("I8").Value = 1.041100180459965
$ws.Range("J8").Value = 1.059896721301164
$ws.Range("K8").Value = 1.064311103120781
$ws.Range("L8").Value = 1.063967263904236
$ws.Range("M8").Value = 1.074214772546179
$ws.Range("N8").Value = 1.023694918248115

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.052416270855899
$ws.Range("D9").Value = 1.059343086414536
$ws.Range("E9").Value = 1.059024404508379
$ws.Range("F9").Value = 1.069152061025528
$ws.Range("I9").Value = 1.040666887208824
$ws.Range("J9").Value = 1.058114627272488
$ws.Range("K9").Value = 1.062431350257568
$ws.Range("L9").Value = 1.062113658790376
$ws.Range("M9").Value = 1.072210154395871
$ws.Range("N9").Value = 1.023090036068289

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.050692344423985
$ws.Range("D10").Value = 1.057807906194489
$ws.Range("E10").Value = 1.057506833699419
$ws.Range("F10").Value = 1.067536834052123
$ws.Range("I10").Value = 1.040370724459911
$ws.Range("J10").Value = 1.056924587669662
$ws.Range("K10").Value = 1.061177355364668
$ws.Range("L10").Value = 1.060877310643747
$ws.Range("M10").Value = 1.070873414197286
$ws.Range("N10").Value = 1.02268548998078

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049946085282966
$ws.Range("D11").Value = 1.057143716743873
$ws.Range("E11").Value = 1.056850323968651
$ws.Range("F11").Value = 1.066838090707149
$ws.Range("I11").Value = 1.040240753264235
$ws.Range("J11").Value = 1.056408819162864
$ws.Range("K11").Value = 1.060634167046795
$ws.Range("L11").Value = 1.060341814773412
$ws.Range("M11").Value = 1.070294514930001
$ws.Range("N11").Value = 1.022510012120831

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049668922689477
$ws.Range("D12").Value = 1.056897090486069
$ws.Range("E12").Value = 1.056606558562936
$ws.Range("F12").Value = 1.066578645650943
$ws.Range("I12").Value = 1.040192216266221
$ws.Range("J12").Value = 1.056217168591504
$ws.Range("K12").Value = 1.060432372643269
$ws.Range("L12").Value = 1.060142885415386
$ws.Range("M12").Value = 1.070079473704228
$ws.Range("N12").Value = 1.022444785821129

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049728373601276
$ws.Range("D13").Value = 1.056949988885685
$ws.Range("E13").Value = 1.056658842925548
$ws.Range("F13").Value = 1.066634292999751
$ws.Range("I13").Value = 1.040202639376557
$ws.Range("J13").Value = 1.056258281511539
$ws.Range("K13").Value = 1.060475659588722
$ws.Range("L13").Value = 1.060185557448791
$ws.Range("M13").Value = 1.070125601331255
$ws.Range("N13").Value = 1.022458779164668

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049923174302076
$ws.Range("D14").Value = 1.057123328841128
$ws.Range("E14").Value = 1.056830172376057
$ws.Range("F14").Value = 1.066816642870986
$ws.Range("I14").Value = 1.040236746482205
$ws.Range("J14").Value = 1.056392978719221
$ws.Range("K14").Value = 1.060617487271951
$ws.Range("L14").Value = 1.060325371666346
$ws.Range("M14").Value = 1.070276739805064
$ws.Range("N14").Value = 1.0225046214327

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.050043201569049
$ws.Range("D15").Value = 1.057230140346287
$ws.Range("E15").Value = 1.056935746238781
$ws.Range("F15").Value = 1.06692900781025
$ws.Range("I15").Value = 1.04025772655881
$ws.Range("J15").Value = 1.056475960680648
$ws.Range("K15").Value = 1.060704868003311
$ws.Range("L15").Value = 1.060411512859748
$ws.Range("M15").Value = 1.070369859583544
$ws.Range("N15").Value = 1.022532860272542

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.050741875606269
$ws.Range("D16").Value = 1.0578519979345
$ws.Range("E16").Value = 1.05755041691819
$ws.Range("F16").Value = 1.067583221272031
$ws.Range("I16").Value = 1.040379313764971
$ws.Range("J16").Value = 1.056958807546783
$ws.Range("K16").Value = 1.061213400766019
$ws.Range("L16").Value = 1.060912846605211
$ws.Range("M16").Value = 1.070911832112813
$ws.Range("N16").Value = 1.022697129412207

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.051180191786536
$ws.Range("D17").Value = 1.058242221077712
$ws.Range("E17").Value = 1.057936146583961
$ws.Range("F17").Value = 1.067993769013603
$ws.Range("I17").Value = 1.040455118914976
$ws.Range("J17").Value = 1.057261557826921
$ws.Range("K17").Value = 1.061532336325555
$ws.Range("L17").Value = 1.06122728015893
$ws.Range("M17").Value = 1.07125177525804
$ws.Range("N17").Value = 1.022800089034501

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.051435874732407
$ws.Range("D18").Value = 1.058469885054718
$ws.Range("E18").Value = 1.058161194869457
$ws.Range("F18").Value = 1.06823329829719
$ws.Range("I18").Value = 1.040499167774075
$ws.Range("J18").Value = 1.05743810125971
$ws.Range("K18").Value = 1.06171834671922
$ws.Range("L18").Value = 1.06141066955602
$ws.Range("M18").Value = 1.071450050433985
$ws.Range("N18").Value = 1.022860114039883

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051523059511551
$ws.Range("D19").Value = 1.058547521645737
$ws.Range("E19").Value = 1.058237940471269
$ws.Range("F19").Value = 1.068314982400284
$ws.Range("I19").Value = 1.040514158966569
$ws.Range("J19").Value = 1.057498290260652
$ws.Range("K19").Value = 1.06178176817886
$ws.Range("L19").Value = 1.061473198162022
$ws.Range("M19").Value = 1.071517655802081
$ws.Range("N19").Value = 1.022880575992539

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.051133162489071
$ws.Range("D20").Value = 1.058200348301292
$ws.Range("E20").Value = 1.057894755368948
$ws.Range("F20").Value = 1.06794971453349
$ws.Range("I20").Value = 1.04044700301273
$ws.Range("J20").Value = 1.057229080312535
$ws.Range("K20").Value = 1.061498119556815
$ws.Range("L20").Value = 1.061193545900931
$ws.Range("M20").Value = 1.071215303393387
$ws.Range("N20").Value = 1.022789045502564

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049865809492561
$ws.Range("D21").Value = 1.057072282245071
$ws.Range("E21").Value = 1.056779717595196
$ws.Range("F21").Value = 1.066762942636062
$ws.Range("I21").Value = 1.040226709963483
$ws.Range("J21").Value = 1.05635331571416
$ws.Range("K21").Value = 1.060575723384393
$ws.Range("L21").Value = 1.060284200480737
$ws.Range("M21").Value = 1.070232233629667
$ws.Range("N21").Value = 1.022491123296845

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049069155022373
$ws.Range("D22").Value = 1.056363504202542
$ws.Range("E22").Value = 1.056079178773072
$ws.Range("F22").Value = 1.066017346798743
$ws.Range("I22").Value = 1.040086699286748
$ws.Range("J22").Value = 1.055802275885934
$ws.Range("K22").Value = 1.059995602290849
$ws.Range("L22").Value = 1.059712329635733
$ws.Range("M22").Value = 1.069614067729065
$ws.Range("N22").Value = 1.022303541622005

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049491459715039
$ws.Range("D23").Value = 1.056739195149362
$ws.Range("E23").Value = 1.056450497407016
$ws.Range("F23").Value = 1.066412546749768
$ws.Range("I23").Value = 1.040161064088138
$ws.Range("J23").Value = 1.05609443157975
$ws.Range("K23").Value = 1.060303151949804
$ws.Range("L23").Value = 1.060015501435247
$ws.Range("M23").Value = 1.069941775875494
$ws.Range("N23").Value = 1.022403007421254

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.051154412940533
$ws.Range("D24").Value = 1.058219268641387
$ws.Range("E24").Value = 1.057913458095874
$ws.Range("F24").Value = 1.067969620659705
$ws.Range("I24").Value = 1.040450670755585
$ws.Range("J24").Value = 1.057243755644331
$ws.Range("K24").Value = 1.061513580701345
$ws.Range("L24").Value = 1.061208789004347
$ws.Range("M24").Value = 1.071231783483428
$ws.Range("N24").Value = 1.022794035690713

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.053084985246027
$ws.Range("D25").Value = 1.059938906674737
$ws.Range("E25").Value = 1.059613444571704
$ws.Range("F25").Value = 1.069779017920592
$ws.Range("I25").Value = 1.040780191530336
$ws.Range("J25").Value = 1.058575690246518
$ws.Range("K25").Value = 1.062917457796283
$ws.Range("L25").Value = 1.062592968703551
$ws.Range("M25").Value = 1.072728454975551
$ws.Range("N25").Value = 1.023246640765683
